# Applies the "Documentation HPA Introduction et analyse" commit:
#  - Fills in two previously-blank journal rows (56 & 57) with a new entry
#    dated 17.05.2022 describing analysis/introduction work on HPA
#    (Hierarchical Pathfinding A*).
#  - Updates the sheet's view state (selected cell / scroll position) to
#    match where the author was working afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 56: "Analyse HPA" session, 08:00 -> 09:30 ---------------------
$ws.Range("A56").Value = 44698                      # 17.05.2022
$ws.Range("B56").Value = 0.33333333333333331        # 08:00
$ws.Range("C56").Value = 0.39583333333333331        # 09:30
$ws.Range("E56").Value = "Analyse HPA"
# D56 already holds the shared "=C56-B56" formula and recalculates on its own.

# --- Row 57: "Introduction et présentation du HPA" session, 09:30 -> 10:27 ---
$ws.Range("A57").Value = 44698                      # 17.05.2022
$ws.Range("B57").Value = 0.39583333333333331        # 09:30
$ws.Range("C57").Value = 0.43541666666666662        # 10:27
$ws.Range("E57").Value = "Introduction et présentation du HPA"
# D57 already holds the shared "=C57-B57" formula and recalculates on its own.

# --- View state: scroll down to keep the new rows in view, select E58 ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 1
$ws.Range("E58").Select()
